$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 11:52"

# Row 4 - Madrid
$ws.Range("B4").Value = 48048
$ws.Range("C4").Value = 26247
$ws.Range("D4").Value = 15233
$ws.Range("E4").Value = 6568

# Row 10 - Navarra
$ws.Range("B10").Value = 4150
$ws.Range("C10").Value = 730
$ws.Range("D10").Value = 3171
$ws.Range("E10").Value = 249

# Row 12 - La Rioja
$ws.Range("B12").Value = 3420
$ws.Range("C12").Value = 1519
$ws.Range("D12").Value = 1677
$ws.Range("E12").Value = 234

# Row 27 - Cantabria
$ws.Range("B27").Value = 1796
$ws.Range("C27").Value = 323
$ws.Range("D27").Value = 1353
$ws.Range("E27").Value = 120
